$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.69843989610672
$ws.Range("B1").Value = 3.397899866104126
$ws.Range("C1").Value = 2.802345514297485
$ws.Range("D1").Value = 2.318350553512573
$ws.Range("E1").Value = 2.078043460845947
